# Update "Performance Metrics" sheet values to match REST API v1.4 formatting:
# remove the trailing zero on "1.50" -> "1.5" and drop the space after the
# comma inside the bracketed confidence interval for both the Odds Ratio (K2)
# and the Concordance Statistic (N2) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance Metrics")

$ws.Range("K2").Value = "1.53 [1.5,1.56]"
$ws.Range("N2").Value = "0.522 [0.519,0.527]"
